$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.313179
$ws.Range("H2").Value = 0.939537
$ws.Range("I2").Value = 0.02707464596575709
$ws.Range("J2").Value = 0.0270746459657571
$ws.Range("M2").Value = 0.5373756666666667
$ws.Range("N2").Value = 1.612127
$ws.Range("O2").Value = 0.007472820128982582
$ws.Range("P2").Value = 0.007472820128982581
$ws.Range("Q2").Value = 0.168294773911
$ws.Range("R2").Value = 1.514652965199
$ws.Range("S2").Value = 0.0002023239593579867
$ws.Range("T2").Value = 0.0002023239593579867

$ws.Range("G3").Value = 0.313179
$ws.Range("H3").Value = 0.939537
$ws.Range("I3").Value = 0.02707464596575709
$ws.Range("J3").Value = 0.0270746459657571
$ws.Range("O3").Value = 0.1537223653287423
$ws.Range("P3").Value = 0.1537223653287423
$ws.Range("Q3").Value = 3.461968878085
$ws.Range("R3").Value = 31.157719902765
$ws.Range("S3").Value = 0.004161978618294472
$ws.Range("T3").Value = 0.004161978618294471

$ws.Range("G4").Value = 0.313179
$ws.Range("H4").Value = 0.939537
$ws.Range("I4").Value = 0.02707464596575709
$ws.Range("J4").Value = 0.0270746459657571
$ws.Range("M4").Value = 30.561198
$ws.Range("N4").Value = 91.683594
$ws.Range("O4").Value = 0.4249882340167162
$ws.Range("P4").Value = 0.4249882340167161
$ws.Range("Q4").Value = 9.571125428442
$ws.Range("R4").Value = 86.14012885597799
$ws.Range("S4").Value = 0.01150640597561492
$ws.Range("T4").Value = 0.01150640597561492

$ws.Range("G5").Value = 0.313179
$ws.Range("H5").Value = 0.939537
$ws.Range("I5").Value = 0.02707464596575709
$ws.Range("J5").Value = 0.0270746459657571
$ws.Range("M5").Value = 29.75783666666667
$ws.Range("N5").Value = 89.27351
$ws.Range("O5").Value = 0.4138165805255589
$ws.Range("P5").Value = 0.4138165805255589
$ws.Range("Q5").Value = 9.31952952943
$ws.Range("R5").Value = 83.87576576487
$ws.Range("S5").Value = 0.01120393741248972
$ws.Range("T5").Value = 0.01120393741248972

$ws.Range("I6").Value = 0.6982806158817221
$ws.Range("J6").Value = 0.6982806158817222
$ws.Range("M6").Value = 0.5373756666666667
$ws.Range("N6").Value = 1.612127
$ws.Range("O6").Value = 0.007472820128982582
$ws.Range("P6").Value = 0.007472820128982581
$ws.Range("Q6").Value = 4.340480703787556
$ws.Range("R6").Value = 39.064326334088
$ws.Range("S6").Value = 0.005218125442039288
$ws.Range("T6").Value = 0.005218125442039288

$ws.Range("I7").Value = 0.6982806158817221
$ws.Range("J7").Value = 0.6982806158817222
$ws.Range("O7").Value = 0.1537223653287423
$ws.Range("P7").Value = 0.1537223653287423
$ws.Range("S7").Value = 0.1073413479365493
$ws.Range("T7").Value = 0.1073413479365493

$ws.Range("I8").Value = 0.6982806158817221
$ws.Range("J8").Value = 0.6982806158817222
$ws.Range("M8").Value = 30.561198
$ws.Range("N8").Value = 91.683594
$ws.Range("O8").Value = 0.4249882340167162
$ws.Range("P8").Value = 0.4249882340167161
$ws.Range("Q8").Value = 246.848338009904
$ws.Range("R8").Value = 2221.635042089136
$ws.Range("S8").Value = 0.296761045791678
$ws.Range("T8").Value = 0.296761045791678

$ws.Range("I9").Value = 0.6982806158817221
$ws.Range("J9").Value = 0.6982806158817222
$ws.Range("M9").Value = 29.75783666666667
$ws.Range("N9").Value = 89.27351
$ws.Range("O9").Value = 0.4138165805255589
$ws.Range("P9").Value = 0.4138165805255589
$ws.Range("Q9").Value = 240.3594428443822
$ws.Range("R9").Value = 2163.23498559944
$ws.Range("S9").Value = 0.2889600967114555
$ws.Range("T9").Value = 0.2889600967114556

$ws.Range("G10").Value = 2.897745666666667
$ws.Range("H10").Value = 8.693237
$ws.Range("I10").Value = 0.2505130868410934
$ws.Range("J10").Value = 0.2505130868410934
$ws.Range("M10").Value = 0.5373756666666667
$ws.Range("N10").Value = 1.612127
$ws.Range("O10").Value = 0.007472820128982582
$ws.Range("P10").Value = 0.007472820128982581
$ws.Range("Q10").Value = 1.557178009455445
$ws.Range("R10").Value = 14.014602085099
$ws.Range("S10").Value = 0.001872039237919684
$ws.Range("T10").Value = 0.001872039237919684

$ws.Range("G11").Value = 2.897745666666667
$ws.Range("H11").Value = 8.693237
$ws.Range("I11").Value = 0.2505130868410934
$ws.Range("J11").Value = 0.2505130868410934
$ws.Range("O11").Value = 0.1537223653287423
$ws.Range("P11").Value = 0.1537223653287423
$ws.Range("Q11").Value = 32.03249679769612
$ws.Range("R11").Value = 288.292471179265
$ws.Range("S11").Value = 0.03850946425501751
$ws.Range("T11").Value = 0.0385094642550175

$ws.Range("G12").Value = 2.897745666666667
$ws.Range("H12").Value = 8.693237
$ws.Range("I12").Value = 0.2505130868410934
$ws.Range("J12").Value = 0.2505130868410934
$ws.Range("M12").Value = 30.561198
$ws.Range("N12").Value = 91.683594
$ws.Range("O12").Value = 0.4249882340167162
$ws.Range("P12").Value = 0.4249882340167161
$ws.Range("Q12").Value = 88.558579072642
$ws.Range("R12").Value = 797.027211653778
$ws.Range("S12").Value = 0.1064651143746725
$ws.Range("T12").Value = 0.1064651143746725

$ws.Range("G13").Value = 2.897745666666667
$ws.Range("H13").Value = 8.693237
$ws.Range("I13").Value = 0.2505130868410934
$ws.Range("J13").Value = 0.2505130868410934
$ws.Range("M13").Value = 29.75783666666667
$ws.Range("N13").Value = 89.27351
$ws.Range("O13").Value = 0.4138165805255589
$ws.Range("P13").Value = 0.4138165805255589
$ws.Range("Q13").Value = 86.23064225020778
$ws.Range("R13").Value = 776.07578025187
$ws.Range("S13").Value = 0.1036664689734836
$ws.Range("T13").Value = 0.1036664689734836

$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.2791366666666666
$ws.Range("H14").Value = 0.83741
$ws.Range("I14").Value = 0.02413165131142748
$ws.Range("J14").Value = 0.02413165131142749
$ws.Range("M14").Value = 0.5373756666666667
$ws.Range("N14").Value = 1.612127
$ws.Range("O14").Value = 0.007472820128982582
$ws.Range("P14").Value = 0.007472820128982581
$ws.Range("Q14").Value = 0.1500012523411111
$ws.Range("R14").Value = 1.35001127107
$ws.Range("S14").Value = 0.0001803314896656242
$ws.Range("T14").Value = 0.0001803314896656243

$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.2791366666666666
$ws.Range("H15").Value = 0.83741
$ws.Range("I15").Value = 0.02413165131142748
$ws.Range("J15").Value = 0.02413165131142749
$ws.Range("O15").Value = 0.1537223653287423
$ws.Range("P15").Value = 0.1537223653287423
$ws.Range("Q15").Value = 3.085655336827778
$ws.Range("R15").Value = 27.77089803145
$ws.Range("S15").Value = 0.00370957451888108
$ws.Range("T15").Value = 0.00370957451888108

$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.2791366666666666
$ws.Range("H16").Value = 0.83741
$ws.Range("I16").Value = 0.02413165131142748
$ws.Range("J16").Value = 0.02413165131142749
$ws.Range("M16").Value = 30.561198
$ws.Range("N16").Value = 91.683594
$ws.Range("O16").Value = 0.4249882340167162
$ws.Range("P16").Value = 0.4249882340167161
$ws.Range("Q16").Value = 8.530750939059999
$ws.Range("R16").Value = 76.77675845154
$ws.Range("S16").Value = 0.01025566787475074
$ws.Range("T16").Value = 0.01025566787475074

$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.2791366666666666
$ws.Range("H17").Value = 0.83741
$ws.Range("I17").Value = 0.02413165131142748
$ws.Range("J17").Value = 0.02413165131142749
$ws.Range("M17").Value = 29.75783666666667
$ws.Range("N17").Value = 89.27351
$ws.Range("O17").Value = 0.4138165805255589
$ws.Range("P17").Value = 0.4138165805255589
$ws.Range("Q17").Value = 8.306503334344443
$ws.Range("R17").Value = 74.7585300091
$ws.Range("S17").Value = 0.009986077428130041
$ws.Range("T17").Value = 0.009986077428130045
